$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.04763786555579896
$ws.Range("C2").Value = 0.04240448674262143
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 0.7443468554461139

# Row 3
$ws.Range("B3").Value = 0.0000000008413942875762359
$ws.Range("C3").Value = 0.00007097389502863649
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 645.32727682996
$ws.Range("G3").Value = 649.227778484905

# Row 4
$ws.Range("B4").Value = 0.0001488876196638067
$ws.Range("C4").Value = 0.00007097389502863649
$ws.Range("D4").Value = 3.900430680208489
$ws.Range("E4").Value = 645.32727682996
$ws.Range("G4").Value = 649.2279273716832
